# Applies the "dashboard auto-update" edit described by the commit:
#  - Re-layout of the Comentarios sheet: a new "extraction_status" column is
#    inserted before the trailing JSON-blob column (which becomes
#    "created_time_raw" shifted one column right, to P).
#  - A FAILED placeholder row is added for every post that was attempted in
#    this run (4 posts total: the 2 that already had comments, plus 2 new
#    ones that failed outright), marking O="FAILED" for those synthetic rows.
#  - The three derived/summary sheets (Resumen_Posts, Stats_Plataforma,
#    Stats_Extraccion) are dropped from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Comentarios")

$lastDataRow = 44

# 1) Shift the trailing JSON-blob column (O) out to a new column P, leaving
#    O free for the new "extraction_status" flag. Also normalize column I's
#    number format onto the same style as the other datetime columns
#    (it used its own yyyy-mm-dd-only style before; now it reuses the
#    shared one).
for ($r = 2; $r -le $lastDataRow; $r++) {
    $jsonVal = $ws.Cells.Item($r, 15).Value2
    $ws.Cells.Item($r, 16).Value = $jsonVal
    $ws.Cells.Item($r, 15).Value = ""
    $ws.Cells.Item($r, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# 2) Header row: O1 becomes "extraction_status", P1 becomes "created_time_raw"
#    (P1 is brand new, so it needs the same bold/bordered header style as
#    its neighbours -- copy it over from O1 before/after the rename).
$ws.Cells.Item(1, 15).Value = "extraction_status"
$ws.Cells.Item(1, 16).Value = "created_time_raw"
$ws.Cells.Item(1, 15).Copy()
$ws.Cells.Item(1, 16).PasteSpecial(-4122)  # xlPasteFormats

# 3) Insert a FAILED placeholder row right after post #1's comments (row 24),
#    pushing the existing post #2 comments down by one row.
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 1
$ws.Cells.Item(24, 2).Value = "Instagram"
$ws.Cells.Item(24, 3).Value = "https://www.instagram.com/p/DSaIkwWgBhr/"
$ws.Cells.Item(24, 4).Value = "https://www.instagram.com/p/DSaIkwWgBhr/"
$ws.Cells.Item(24, 5).Value = ""
$ws.Cells.Item(24, 6).Value = ""
$ws.Cells.Item(24, 7).Value = ""
$ws.Cells.Item(24, 8).Value = ""
$ws.Cells.Item(24, 9).Value = ""
$ws.Cells.Item(24, 10).Value = ""
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = $false
$ws.Cells.Item(24, 14).Value = ""
$ws.Cells.Item(24, 15).Value = "FAILED"
$ws.Cells.Item(24, 16).Value = ""
# The row-insert above copied row 23's date format onto H24/I24; this is a
# brand-new placeholder row with no date in it, so drop back to the
# workbook's default (unstyled) cell, matching the other blank cells here.
$ws.Cells.Item(24, 8).ClearFormats()
$ws.Cells.Item(24, 9).ClearFormats()

# 4) Append three more FAILED placeholder rows at the bottom of the sheet,
#    one per post that failed to extract in this run (post #2 again, plus
#    two brand-new advertiser posts that never produced any comments).
$failedPosts = @(
    @{ Num = 2; Url = "https://www.instagram.com/p/DRiT7h1gO2m/" },
    @{ Num = 3; Url = "https://www.instagram.com/p/DSnPl7mgFk-/#advertiser" },
    @{ Num = 4; Url = "https://www.instagram.com/p/DSnPlKqgHU7/#advertiser" }
)

$row = 46
foreach ($post in $failedPosts) {
    $ws.Cells.Item($row, 1).Value = $post.Num
    $ws.Cells.Item($row, 2).Value = "Instagram"
    $ws.Cells.Item($row, 3).Value = $post.Url
    $ws.Cells.Item($row, 4).Value = $post.Url
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""
    $ws.Cells.Item($row, 9).Value = ""
    $ws.Cells.Item($row, 10).Value = ""
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = $false
    $ws.Cells.Item($row, 14).Value = ""
    $ws.Cells.Item($row, 15).Value = "FAILED"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 8).ClearFormats()
    $ws.Cells.Item($row, 9).ClearFormats()
    $row++
}

# 5) Drop the derived/summary sheets; only the raw comments sheet remains.
$wb.Worksheets("Resumen_Posts").Delete()
$wb.Worksheets("Stats_Plataforma").Delete()
$wb.Worksheets("Stats_Extraccion").Delete()
